$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new GitHub admin-log entry as row 19.
$ws.Range("A19").Value = "2025-09-03 06:54:24"
$ws.Range("B19").Value = "create-team"
$ws.Range("C19").Value = "new-organization97"
$ws.Range("D19").Value = "newtestteam"
$ws.Range("E19").Value = "demo"

# "False" would otherwise be auto-recognised as the Boolean FALSE by
# Excel's type inference; force it to remain literal text (matching the
# other rows in column I), then drop the resulting quote-prefix style so
# no extraneous cell format is introduced.
$ws.Range("I19").Value = "'False"
$ws.Range("I19").Style = "Normal"
